$wb = $excel.ActiveWorkbook

# New phone-number values to drop into the "MobileNumber" / "*_PN" columns
# (G2, AF2, AV2, AZ2) of the row-2 sample record on each of the first four
# sheets. These must land as plain TEXT (shared-string) cells, matching the
# pre-existing cells' type/style ("General" style, s=2) rather than being
# auto-coerced to numbers by a naive Value assignment.
$updates = [ordered]@{
    "G2"  = "9840099433"
    "AF2" = "9840038682"
    "AV2" = "9840096039"
    "AZ2" = "9840036951"
}

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    foreach ($addr in $updates.Keys) {
        $cell = $ws.Range($addr)
        # Writing the digit string straight through .Value coerces it to a
        # number. Route it through a text-returning formula first, then
        # collapse the formula down to its literal text result (Copy +
        # PasteSpecial values-only) so the stored cell keeps its original
        # style and becomes a genuine string cell, just like the source
        # workbook's existing phone-number entries.
        $cell.Formula = '="' + $updates[$addr] + '"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)
    }
}

$excel.CutCopyMode = $false
